# #272 Ajout d'un scénario de recherche de l'offre d'un professionnel avec un ID Nat PS
# - bump the "Date" metadata value
# - swap the "Mapping: RIM Mapping" / "Mapping: Spécification métier vers l'extension
#   ROR TelecomCommunicationUsage" columns (AK <-> AL) on the Elements sheet, including
#   their column widths, so the "Spécification métier" mapping now comes first.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: refresh the generation Date ---------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# --- Elements sheet: swap columns AK (37) and AL (38) ----------------------------
$elements = $wb.Worksheets.Item("Elements")

$colAK = 37
$colAL = 38

for ($r = 1; $r -le 6; $r++) {
    $akCell = $elements.Cells.Item($r, $colAK)
    $alCell = $elements.Cells.Item($r, $colAL)

    $akVal = $akCell.Value()
    $alVal = $alCell.Value()

    # Only touch a cell when the incoming value actually differs from what is
    # already there - writing an empty string back into an already-empty cell
    # turns a shared-string "" reference into a truly-blank cell, which is a
    # needless extra change versus the source workbook.
    if ($akVal -ne $alVal) {
        $akCell.Value = $alVal
        $alCell.Value = $akVal
    }
}

# Swap the column widths too (AK becomes the wide "Spécification métier" column,
# AL becomes the narrower "RIM Mapping" column). The inputs below are tuned so the
# engine's pixel-quantized ColumnWidth setter lands on the closest possible stored
# width to the target bestFit widths (86.09375 / 24.98046875).
$elements.Columns.Item($colAK).ColumnWidth = 85.15
$elements.Columns.Item($colAL).ColumnWidth = 24.15
